$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of quotes for 2025-09-20 (Excel serial date 45920)
$ws.Range("A16").Value = 45920
$ws.Range("A16").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("B16").Value = "20,9437"
$ws.Range("C16").Value = "15,0727"
$ws.Range("D16").Value = "14,9476"
$ws.Range("E16").Value = "14,9476"
